# PM04 Tidsregistrering for Toke.xlsx - add new time-tracking entries (rows 13-21)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 13: UC01 og UC02 frem til SSD
$ws.Range("A13").Value = "UC01 og UC02 frem til SSD"
$ws.Range("C13").Value = 43964
$ws.Range("D13").Value = 0.354166666666667
$ws.Range("E13").Value = 0.583333333333333

# Row 14: Vejledning fra Andrés
$ws.Range("A14").Value = "Vejledning fra Andrés"
$ws.Range("C14").Value = 43964
$ws.Range("D14").Value = 0.583333333333333
$ws.Range("E14").Value = 0.666666666666667

# Row 15: Vejledning fra Anders
$ws.Range("A15").Value = "Vejledning fra Anders"
$ws.Range("C15").Value = 43964
$ws.Range("D15").Value = 0.666666666666667
$ws.Range("E15").Value = 0.729166666666667

# Row 16: Rettelse af UCD01
$ws.Range("A16").Value = "Rettelse af UCD01"
$ws.Range("C16").Value = 43965
$ws.Range("D16").Value = 0.354166666666667
$ws.Range("E16").Value = 0.40625

# Row 17: Masse omdøbelse af filer
$ws.Range("A17").Value = "Masse omdøbelse af filer"
$ws.Range("C17").Value = 43965
$ws.Range("D17").Value = 0.40625
$ws.Range("E17").Value = 0.427083333333333

# Row 18: Rettelse af UC02, samt vejledning fa Anders
$ws.Range("A18").Value = "Rettelse af UC02, samt vejledning fa Anders"
$ws.Range("C18").Value = 43965
$ws.Range("D18").Value = 0.427083333333333
$ws.Range("E18").Value = 0.572916666666667

# Row 19: Rettelse af AD01
$ws.Range("A19").Value = "Rettelse af AD01"
$ws.Range("C19").Value = 43965
$ws.Range("D19").Value = 0.572916666666667
$ws.Range("E19").Value = 0.604166666666667

# Row 20: Rettelse af DOM02
$ws.Range("A20").Value = "Rettelse af DOM02"
$ws.Range("C20").Value = 43965
$ws.Range("D20").Value = 0.604166666666667
$ws.Range("E20").Value = 0.645833333333333

# Row 21: Rettelse af SSD01
$ws.Range("A21").Value = "Rettelse af SSD01"
$ws.Range("C21").Value = 43965
$ws.Range("D21").Value = 0.645833333333333
$ws.Range("E21").Value = 0.6875

# Update the selected cell in the sheet view to E22
$ws.Range("E22").Select()
